$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: header + text number format
$ws.Range("K1").Value = "Дата для решения"
$ws.Columns("K").ColumnWidth = 15.333333333333332

# Apply text number format ("@") to columns B, D, E, K (new style xf: numFmtId 49)
$ws.Range("B1").NumberFormat = "@"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("K1").NumberFormat = "@"

# Clear the example data row (row 2), keep H2's hyperlink style but remove its value
$ws.Range("A2:G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2:K2").ClearContents()

# Selection / view state
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("F11").Select()
